# Apply updated GHI weather data values to the workbook.
$wb = $excel.ActiveWorkbook

$daily = $wb.Worksheets.Item("Daily")
$hourly = $wb.Worksheets.Item("Hourly")

# --- Daily sheet, row 2 ---
$daily.Range("G2").Value = 4073.79
$daily.Range("H2").Value = 7403.28
$daily.Range("I2").Value = 883.05
$daily.Range("J2").Value = 1972.7
$daily.Range("K2").Value = 459.33
$daily.Range("L2").Value = 1712.5

# --- Hourly sheet ---

# row 9
$hourly.Range("K9").Value = 22.78
$hourly.Range("M9").Value = 22.78

# row 10
$hourly.Range("K10").Value = 68.36
$hourly.Range("M10").Value = 68.36

# row 11
$hourly.Range("K11").Value = 118.79
$hourly.Range("L11").Value = 0
$hourly.Range("M11").Value = 118.79

# row 12
$hourly.Range("K12").Value = 199.9
$hourly.Range("L12").Value = 14.35
$hourly.Range("M12").Value = 192.41

# row 13
$hourly.Range("I13").Value = 833.46
$hourly.Range("K13").Value = 285.03
$hourly.Range("L13").Value = 37.26
$hourly.Range("M13").Value = 263.61

# row 14
$hourly.Range("H14").Value = 604.98
$hourly.Range("I14").Value = 845.49
$hourly.Range("J14").Value = 104.46
$hourly.Range("K14").Value = 340.34
$hourly.Range("L14").Value = 107.4
$hourly.Range("M14").Value = 275.86

# row 15
$hourly.Range("H15").Value = 579.9400000000001
$hourly.Range("I15").Value = 835.6
$hourly.Range("K15").Value = 377.51
$hourly.Range("L15").Value = 183.32
$hourly.Range("M15").Value = 271.34

# row 16
$hourly.Range("H16").Value = 502.27
$hourly.Range("I16").Value = 801.55
$hourly.Range("J16").Value = 96.78
$hourly.Range("K16").Value = 291.5
$hourly.Range("L16").Value = 108.86
$hourly.Range("M16").Value = 234.1

# row 17
$hourly.Range("H17").Value = 380.09
$hourly.Range("I17").Value = 734.47
$hourly.Range("K17").Value = 149.84
$hourly.Range("L17").Value = 6.9
$hourly.Range("M17").Value = 147

# row 18
$hourly.Range("H18").Value = 227.86
$hourly.Range("I18").Value = 609.46
$hourly.Range("K18").Value = 91.45
$hourly.Range("L18").Value = 1.23
$hourly.Range("M18").Value = 91.05

# row 19
$hourly.Range("I19").Value = 335.53
$hourly.Range("K19").Value = 27.19
$hourly.Range("M19").Value = 27.19

$wb.Save()
